# "Aankomst en Vertrek Depots toegevoegd"
# Append the arrival/departure depot edges to the edge list on Sheet1
# (columns: A=from, B=to, C=length), then leave the selection where the
# author left it after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each triple is (from, to, length) - mirrors the existing A:C layout.
$newEdges = @(
    @(29, 109, 0.5),
    @(109, 29, 0.5),
    @(109, 112, 0.5),
    @(112, 109, 0.5),
    @(104, 110, 0.5),
    @(110, 104, 0.5),
    @(110, 112, 0.5),
    @(112, 110, 0.5),
    @(10, 111, 0.5),
    @(111, 10, 0.5),
    @(111, 113, 0.5),
    @(113, 111, 0.5)
)

$xlUp = -4162
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row + 1

for ($i = 0; $i -lt $newEdges.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value2 = $newEdges[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $newEdges[$i][1]
    $ws.Cells.Item($row, 3).Value2 = $newEdges[$i][2]
}

$ws.Range("E262").Select()
